$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 706
$endRow = 713

$dates = @("07/06/2023 23:56","07/06/2023 23:56","07/06/2023 23:56","07/06/2023 23:56","07/06/2023 23:56","07/06/2023 23:56","07/06/2023 23:56","07/06/2023 23:56")
$coins = @("Germania","American Eagle","Panda","Wedge Tailed","Koala","Krugerrand","Britannia","Filarmonica")
$prices = @("1927.2","36.01","34.73","1927.2","229.31","28.81","229.31","27.76")
$stocks = @("`t`tFuera de Stock`t","`t`tDisponible  ","`t`tDisponible  ","`t`tAgotado Temporalmente`t","`t`tFuera de Stock`t","`t`tDisponible  ","`t`tFuera de Stock`t","`t`tDisponible  ")

# Force column C (Precio) to be stored as text, matching the rest of the sheet,
# otherwise Excel would auto-detect these numeric-looking strings as numbers.
$ws.Range(("C{0}:C{1}" -f $startRow, $endRow)).NumberFormat = "@"

for ($i = 0; $i -lt ($endRow - $startRow + 1); $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $coins[$i]
    $ws.Cells.Item($row, 3).Value = $prices[$i]
    $ws.Cells.Item($row, 4).Value = $stocks[$i]
}
